$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

$data1 = @(
  @("NEI-CEDA CI", 0, 4, 3695, 895, "🟡 Observer", "➖ Neutre"),
  @("BRVM - SERVICES PUBLICS", 0, 8, 3336.84, 111.31, "🟡 Observer", "➖ Neutre"),
  @("BRVM - AUTRES SECTEURS", 0, 4, 2379.5, 596.54, "🟡 Observer", "➖ Neutre"),
  @("BRVM - DISTRIBUTION", 0, 4, 2009.95, 502.65, "🟡 Observer", "➖ Neutre"),
  @("BRVM - TRANSPORT", 0, 4, 1436.81, 362.26, "🟡 Observer", "➖ Neutre"),
  @("BRVM - AGRICULTURE", 0, 4, 1348.69, 336.05, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 683.05, 170.7, "🟡 Observer", "➖ Neutre"),
  @("BRVM - FINANCES", 0, 4, 588.37, 147.3, "🟡 Observer", "➖ Neutre"),
  @("BRVM - SERVICES FINANCIERS", 0, 4, 578.25, 144.77, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRESTIGE", 0, 4, 566.17, 141.59, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIELS", 0, 4, 509.37, 128.91, "🟡 Observer", "➖ Neutre"),
  @("BRVM - ENERGIE", 0, 4, 449.89, 113, "🟡 Observer", "➖ Neutre"),
  @("BRVM - TELECOMMUNICATIONS", 0, 4, 376.2, 93.68000000000001, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE                 (**)", 0, 1, 266.86, 266.86, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE", 0, 1, 266.4, 266.4, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIE    (**)", 0, 1, 262.27, 262.27, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE         (**)", 0, 1, 222.7, 222.7, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE", 0, 1, 222.06, 222.06, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL                    (**)", 0, 1, 220.55, 220.55, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL", 0, 1, 220.02, 220.02, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRINCIPAL     (**)", 0, 1, 219.45, 219.45, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE   (**)", 0, 1, 218.71, 218.71, "🟡 Observer", "➖ Neutre"),
  @("ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 4, 0, 29.37, 7.32, "🟢 Achat", "✅ Renforcer"),
  @("SUCRIVOIRE (SCRC)", 2, 0, 4.62, 2.76, "🟡 Observer", "➖ Neutre"),
  @("ECOBANK TRANS. INCORP. TG (ETIT)", 1, 0, 4.55, 4.55, "🟡 Observer", "➖ Neutre"),
  @("FILTISAC CI (FTSC)", 1, 0, 2.91, 2.91, "🟡 Observer", "➖ Neutre"),
  @("ECOBANK COTE D''IVOIRE (ECOC)", 1, 0, 2.89, 2.89, "🟡 Observer", "➖ Neutre"),
  @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 1, 0, 2.76, 2.76, "🟡 Observer", "➖ Neutre"),
  @("TOTALENERGIES MARKETING SN (TTLS)", 1, 0, 2.41, 2.41, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA NG (BOAN)", 1, 1, 2.21, -1.35, "🟡 Observer", "👀 À surveiller"),
  @("SMB CI (SMBC)", 1, 1, 2.14, 3.19, "🟡 Observer", "👀 À surveiller"),
  @("CORIS BANK INTERNATIONAL (CBIBF)", 1, 0, 1.6, 1.6, "🟡 Observer", "➖ Neutre"),
  @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
  @("ONATEL BF (ONTBF)", 2, 1, -0.1, 2.04, "🟡 Observer", "👀 À surveiller"),
  @("SERVAIR ABIDJAN CI (ABJC)", 1, 1, -0.8100000000000001, -4.73, "🟡 Observer", "👀 À surveiller"),
  @("BICI CI (BICC)", 0, 1, -1.64, -1.64, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -1.76, -1.76, "🟡 Observer", "➖ Neutre"),
  @("SICOR CI (SICC)", 1, 2, -2.23, -3.71, "🟡 Observer", "👀 À surveiller"),
  @("BERNABE CI (BNBC)", 2, 1, -2.28, -1.38, "🟡 Observer", "👀 À surveiller"),
  @("BANK OF AFRICA CI (BOAC)", 0, 1, -2.3, -2.3, "🟡 Observer", "➖ Neutre"),
  @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 0, 1, -2.42, -2.42, "🟡 Observer", "➖ Neutre"),
  @("ORAGROUP TOGO (ORGT)", 0, 1, -2.72, -2.72, "🟡 Observer", "➖ Neutre"),
  @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 0, 1, -3.48, -3.48, "🟡 Observer", "➖ Neutre"),
  @("VIVO ENERGY CI (SHEC)", 0, 1, -3.85, -3.85, "🟡 Observer", "➖ Neutre"),
  @("SICABLE CI (CABC)", 0, 1, -3.89, -3.89, "🟡 Observer", "➖ Neutre"),
  @("NEI-CEDA CI (NEIC)", 0, 1, -4.26, -4.26, "🟡 Observer", "➖ Neutre"),
  @("ORANGE COTE D'IVOIRE (ORAC)", 0, 2, -5.98, -4.03, "🟡 Observer", "➖ Neutre"),
  @("UNILEVER CI (UNLC)", 0, 1, -7.5, -7.5, "🟡 Observer", "➖ Neutre")
)

for ($i = 0; $i -lt $data1.Count; $i++) {
    $row = $data1[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws1.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

$ws1.Range("A50:G53").EntireRow.Delete()

$data2 = @(
  @("BRVM - SERVICES PUBLICS", 9151742.550000001),
  @("NEI-CEDA CI", 1097744.94),
  @("BRVM - AUTRES SECTEURS", 233042.66),
  @("BRVM - DISTRIBUTION", 131662.57),
  @("BRVM - TRANSPORT", 44356.82),
  @("BRVM - AGRICULTURE", 36426.47),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 5274.69),
  @("BRVM - FINANCES", 3627.66),
  @("BRVM - SERVICES FINANCIERS", 3477.31),
  @("BRVM-PRESTIGE", 3303.83)
)

for ($i = 0; $i -lt $data2.Count; $i++) {
    $row = $data2[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws2.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
